$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address values in column A (rows 2 and 3):
# old "sonar.mayuresh260197@gmail.com" -> new "kumarbob912@gmail.com"
$ws.Range("A2").Value = "kumarbob912@gmail.com"
$ws.Range("A3").Value = "kumarbob912@gmail.com"

# Update the hyperlink attached to A3 to point at the new email address
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$3') {
        $h.Address = "mailto:kumarbob912@gmail.com"
    }
}

# Remove the hyperlink that used to live on A2 (profile page no longer links A2)
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.Delete()
    }
}

# Move the active selection to A3
$ws.Range("A3").Select()
